# Add the I0 (I) and IF (J) columns to the sheet, as described in the
# commit "I0 and IF added".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (bold font + border) of the existing "IP" header cell
# onto the two new header cells before we set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Header row: I1 = "I0", J1 = "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-47.
$iValues = @(1,1,1,6,1,1,1,1,1,1,1,1,1,5,1,1,1,1,1,8,5,5,3,5,8,6,7,6,8,6,8,6,8,6,7,7,9,7,3,5,6,5,2,7,7,5)
$jValues = @(5,5,5,6,5,4,7,6,6,6,6,5,3,7,5,5,3,3,3,8,6,6,4,6,8,7,7,7,8,6,8,6,8,6,7,7,9,8,4,5,6,6,2,7,7,5)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
